$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header row grows from B1 to E1, and B1's label changes ---
$ws.Range("B1").Value = "truc"
$ws.Range("C1").Value = "chose"
$ws.Range("D1").Value = "machin"
$ws.Range("E1").Value = "autre"

# --- Row 2 ---
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = 1.95
$ws.Range("D2").Value = 1.45
$ws.Range("E2").Value = 2.25

# --- Row 3 ---
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = 1.78
$ws.Range("E3").Value = 4

# --- Row 4 ---
$ws.Range("C4").Value = 0.28
$ws.Range("D4").Value = 2.5
$ws.Range("E4").Value = 5

# --- Row 5 ---
$ws.Range("C5").Value = 2.02
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = 0.75

# --- Row 6: new "orange" label in E6 ---
$ws.Range("E6").Value = "orange"

# --- Selection moves to E2 ---
$ws.Range("E2").Select()
